$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.799.56"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "3.711.82"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'597.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").Value = "'166.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("D7").Value = "3.707.42"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  +4.43%  "
$ws.Range("D11").Value = "'6.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").Value = "'0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").Value = "'38.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "4.334.17"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "3.715.91"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").Value = "67.780.64"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'7.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").Value = "'17.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.59%  "
$ws.Range("D21").Value = "'488.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "'9.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "'0.726"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("D24").Value = "'84.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("E25").Value = "  +3.78%  "
$ws.Range("D26").Value = "'2.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("D28").Value = "'10.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'7.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.45%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.44%  "
$ws.Range("D33").Value = "'31.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.00%  "
$ws.Range("D34").Value = "3.854.39"
$ws.Range("E34").Value = "  -1.75%  "
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("D36").Value = "3.657.83"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").Value = "'0.132"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("D41").Value = "'0.323"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'48.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'428.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("D45").Value = "'2.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").Value = "'8.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "'40.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").Value = "'141.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.83%  "
$ws.Range("D50").Value = "'0.0352"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").Value = "2.755.47"
$ws.Range("E51").Value = "  -3.03%  "
